# Remove the data rows for the deleted job-application entries (ids 1-3,
# previously on rows 3-5), shifting everything below them up - this is the
# effect of the new "delete row" function.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3:E5").EntireRow.Delete() | Out-Null

# The remaining entry (row 2) gets its applyDate refreshed to a new
# ISO-8601 timestamp string.
$ws.Range("E2").Value = "2021-08-12T00:01:21.334Z"
